$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3815
$ws.Range("A3").Value = 3759
$ws.Range("C3").Value = 2
$ws.Range("A4").Value = 3247
$ws.Range("C4").Value = 4
$ws.Range("A5").Value = 8918
$ws.Range("C5").Value = 5
$ws.Range("A6").Value = 2456
$ws.Range("C6").Value = 3

$wb.Save()
